# Chapter 2 slide: turn the trailing word "here" (in the *SBManager helper
# class* caption) into a hyperlink pointing at the class's source code.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item("TextBox 2")
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$idx = $full.LastIndexOf("here", [StringComparison]::Ordinal)
$startPos = $idx + 1
$len = "here".Length

$sub = $tr.Characters($startPos, $len)
$sub.ActionSettings(1).Hyperlink.Address = "https://github.com/Programming5554/2024Robot/blob/main/src/main/java/frc/robot/SBManager.java"
